$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 421, shifting existing rows 421..466 down to 422..467
$ws.Rows.Item(421).Insert()

# Populate the new row 421 with the new data record
$ws.Cells.Item(421, 1).Value = 7
$ws.Cells.Item(421, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(421, 3).Value = "Ñuble"
$ws.Cells.Item(421, 4).Value = 44918
$ws.Cells.Item(421, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(421, 5).Value = 16
$ws.Cells.Item(421, 6).Value = 100114001
$ws.Cells.Item(421, 7).Value = "Papa"
$ws.Cells.Item(421, 8).Value = "Asterix"
$ws.Cells.Item(421, 9).Value = "1a nueva(o)"
$ws.Cells.Item(421, 10).Value = 200
$ws.Cells.Item(421, 11).Value = 12000
$ws.Cells.Item(421, 12).Value = 12000
$ws.Cells.Item(421, 13).Value = 12000
$ws.Cells.Item(421, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(421, 15).Value = "Región del Maule"
$ws.Cells.Item(421, 16).Value = 480
$ws.Cells.Item(421, 17).Value = 25
$ws.Cells.Item(421, 18).Value = "Hortaliza"
